$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 10068.219999999999
$ws.Range("B10").Value = 10125.94
$ws.Range("C10").Value = 307.87
$ws.Range("D10").Value = 306.13
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = -0.56999999999999995
$ws.Range("G10").Value = 42612.67292824074
$ws.Range("G10").NumberFormat = "m/d/yy h:mm"
$ws.Range("H10").Value = $false
